$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 17.7159309387207
$ws.Range("C3").Value = 17.00711250305176
$ws.Range("C4").Value = 16.59107208251953
$ws.Range("C5").Value = 16.47686958312988
$ws.Range("C6").Value = 16.36934280395508
